$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Delete rows 3 through 23 (everything below the header + first data row),
# leaving only the header row and one data row.
$ws.Range("A3:B23").EntireRow.Delete() | Out-Null

# Update the remaining data row (row 2) to hold the 2008 entry.
$ws.Range("A2").Value = "C:\Users\zaka\Desktop\MOTOGP\Excels\data\2008.xlsx"
$ws.Range("B2").Value = "2008"

# Match the selection recorded in the saved file.
$ws.Range("B2").Select() | Out-Null
